$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (admin_level_2 is D, collector_name was E)
# so a new "village" column lands right after "admin_level_2".
$ws.Range("E1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("E1").Value2 = "village"
